$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 1: translate existing group headers and add a new "Di truyen"
# (genetic algorithm) group between the "Evolutionary" and "Manual"
# columns, shifting "Manual" from L1 to Q1.
# ------------------------------------------------------------------

# Re-create the merges so the brand new L1:P1 range is inserted first
# (matching the target merge order), then restore the original two
# header merges.
$ws.Cells.UnMerge()
$ws.Range("L1:P1").Merge()
$ws.Range("B1:F1").Merge()
$ws.Range("G1:K1").Merge()

# Merging re-styles the affected ranges (Excel assigns them a fresh
# cell style with auto borders). Re-apply the original header style
# (same as A1, which was untouched) to the whole header row so every
# cell keeps/returns to the shared style used by the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("B1:Q1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "Tham lam"
$ws.Range("G1").Value = "Quy hoạch động"
$ws.Range("L1").Value = "Di truyền"
$ws.Range("Q1").Value = "Manual"

# ------------------------------------------------------------------
# Row 2: statistic labels. Add mean/std/min/med/max for the new L:P
# block and move "optim." from L2 to Q2.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("L2:Q2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("L2").Value = "mean"
$ws.Range("M2").Value = "std"
$ws.Range("N2").Value = "min"
$ws.Range("O2").Value = "med"
$ws.Range("P2").Value = "max"
$ws.Range("Q2").Value = "optim."

# ------------------------------------------------------------------
# Row 4 (Problem 1): extend with the new genetic-algorithm block and
# shift the "optim." value into the new Q column.
# ------------------------------------------------------------------
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 50
$ws.Range("O4").Value = 100
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 100

# ------------------------------------------------------------------
# Row 5 (Problem 2): updated values, extended with the new
# genetic-algorithm block.
# ------------------------------------------------------------------
$ws.Range("B5").Value = 36
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 35
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 45
$ws.Range("G5").Value = 35
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 35
$ws.Range("L5").Value = 45
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 45
$ws.Range("O5").Value = 45
$ws.Range("P5").Value = 45
$ws.Range("Q5").Value = 20

# ------------------------------------------------------------------
# New rows 6-13: Problem 3 .. Problem 10
# ------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A6:A13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$data = @(
    @("Problem 3",  171, 29.14, 140, 175, 200, 200, 0, 200, 200, 200, 200, 0, 200, 200, 200, 60),
    @("Problem 4",  200, 0,     200, 200, 200, 200, 0, 200, 200, 200, 200, 0, 200, 200, 200, 160),
    @("Problem 5",  175, 0,     175, 175, 175, 175, 0, 175, 175, 175, 175, 0, 175, 175, 175, 105),
    @("Problem 6",  235, 0,     235, 235, 235, 235, 0, 235, 235, 235, 235, 0, 235, 235, 235, 140),
    @("Problem 7",  300, 0,     300, 300, 300, 300, 0, 300, 300, 300, 300, 0, 300, 300, 300, 300),
    @("Problem 8",  225, 0,     225, 225, 225, 225, 0, 225, 225, 225, 225, 0, 225, 225, 225, 135),
    @("Problem 9",  275, 0,     275, 275, 275, 275, 0, 275, 275, 275, 275, 0, 275, 275, 275, 275),
    @("Problem 10", 300, 0,     300, 300, 300, 300, 0, 300, 300, 300, 300, 0, 300, 300, 300, 300)
)

$rowIndex = 6
foreach ($rowData in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowData[$col - 1]
    }
    $rowIndex++
}
